# Insert a new "Our Brief" slide right after the title slide (slide 1),
# using the same "Title and Content" layout as the rest of the deck.
$p = $ppt.ActivePresentation

$ppLayoutText = 2
$newSlide = $p.Slides.Add(2, $ppLayoutText)

$titleRange = $newSlide.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "Our Brief"
$titleRange.LanguageID = "en-GB"
